$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply corrected values for Temps_Comp_Descompensacio (H) and Comp_Descompensacio (G)
# columns. These reflect fixing the Competing risks calculation: when a
# decompensation event occurred (Descompensacio/Comp_Descompensacio = 1) the
# time value must be taken from Temps_Descompensacio (col B) instead of 0; for the
# two rows where a death also occurred simultaneously (rows 250 and 308), the
# competing-event code is set to 2 and the time taken from Temps_Mort (col F).
$updates = @{
    "H3" = 48
    "H12" = 4
    "H13" = 40
    "H15" = 21
    "H16" = 12
    "H17" = 13
    "H18" = 4
    "H22" = 23
    "H25" = 12
    "H27" = 120
    "H28" = 12
    "H33" = 126.8
    "H39" = 9
    "H41" = 6
    "H47" = 54
    "H48" = 11
    "H49" = 8
    "H53" = 10
    "H57" = 26
    "H60" = 9
    "H67" = 61
    "H69" = 12
    "H71" = 156
    "H81" = 2
    "H89" = 32
    "H91" = 35
    "H101" = 12
    "H104" = 5
    "H105" = 16
    "H114" = 23
    "H119" = 19
    "H127" = 33
    "H133" = 68
    "H146" = 28
    "H149" = 13
    "H153" = 7
    "H155" = 10
    "H156" = 61
    "H160" = 39
    "H166" = 3
    "H167" = 2
    "H172" = 76
    "H177" = 2
    "H178" = 31
    "H180" = 37
    "H193" = 20
    "H199" = 1
    "H220" = 1
    "H225" = 16
    "H229" = 8
    "H231" = 31
    "H232" = 1
    "H236" = 1
    "H240" = 9
    "H246" = 1
    "G250" = 2
    "H250" = 1
    "H269" = 1
    "H285" = 1
    "H287" = 1
    "H290" = 45
    "G308" = 2
    "H308" = 31
    "H318" = 1
    "H319" = 36
    "H322" = 1
    "H328" = 1
    "H329" = 86
    "H336" = 23
    "H342" = 76
    "H354" = 1
    "H361" = 1
    "H364" = 1
    "H371" = 4
    "H372" = 1
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
